$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("B2")
$r.Value = "再形成"
$r.SetPhonetic()
Write-Output ("Count after SetPhonetic: " + $r.Phonetics.Count)
try {
  $ph = $r.Phonetics
  Write-Output ("Phonetics type acquired")
} catch {
  Write-Output ("error: " + $_)
}
